# feat: add 2022-Q1 data
#
# - Repurpose the old "总计" (sheet index 3) worksheet as the new "2022-Q1"
#   worksheet: rename it and fill it in with the new quarter's fund-holdings
#   table (same shape as the "2021-Q3"/"2021-Q4" sheets).
# - Append a brand-new "总计" worksheet after it (cloned from "2022-Q1" so it
#   keeps the same sheet/page setup), holding the previous summary rows plus
#   a new leading row for "2022-Q1".

$wb = $excel.ActiveWorkbook

# Helper: write a value as genuine TEXT (keeps leading zeros / trailing
# zeros, e.g. "002861" or "0.0590") and then reset the cell style back to
# "Normal" so no stray number-format is left applied to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Rename the existing "总计" sheet (3rd tab) to "2022-Q1" and replace
#    its contents with the new fund-holdings table.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Pull header/index-column formatting from the "2021-Q4" sheet so the new
# sheet matches the look of its siblings (bold header row + bordered A col).
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

# Headers
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "002861"
Set-TextValue $q1.Range("C2") "工银瑞信智能制造股票"
Set-TextValue $q1.Range("D2") "1.34"
Set-TextValue $q1.Range("E2") "93.93"
Set-TextValue $q1.Range("F2") "4.40"
Set-TextValue $q1.Range("G2") "0.0590"
$q1.Range("H2").Value = 9

# Row 3
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "001721"
Set-TextValue $q1.Range("C3") "工银瑞信新增益混合"
Set-TextValue $q1.Range("D3") "5.57"
Set-TextValue $q1.Range("E3") "30.21"
Set-TextValue $q1.Range("F3") "1.03"
Set-TextValue $q1.Range("G3") "0.0574"
$q1.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2) Add a new "总计" sheet after "2022-Q1" (cloned from it, then wiped and
#    rebuilt) with the updated summary: previous two rows, plus a new
#    leading "2022-Q1" row.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Copy($null, $lastSheet)
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计"

# Wipe the cloned fund-holdings table and rebuild as the summary table.
$total.Cells.Clear()

$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.12

$total.Range("A3").Value = 1
Set-TextValue $total.Range("B3") "2021-Q4"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 0.62

$total.Range("A4").Value = 2
Set-TextValue $total.Range("B4") "2021-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.08

# Restore the originally-active sheet/selection.
$wb.Worksheets.Item(1).Activate()
